$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = -7.728999999999999
$ws.Range("C3").Value = -11.107
$ws.Range("D3").Value = -7.303
$ws.Range("C4").Value = -12.943
$ws.Range("E8").Value = 16.746
$ws.Range("D9").Value = -7.223000000000001
$ws.Range("A11").Value = -21.794
$ws.Range("E11").Value = 17.182
$ws.Range("A12").Value = -21.683
$ws.Range("C14").Value = -12.592
$ws.Range("E14").Value = 17.256
$ws.Range("A15").Value = -21.791
$ws.Range("D15").Value = -7.813
$ws.Range("E15").Value = 16.741
$ws.Range("E17").Value = 16.787
$ws.Range("D19").Value = -8.103
$ws.Range("D20").Value = -7.825
$ws.Range("D25").Value = -7.802
$ws.Range("C26").Value = -12.214
$ws.Range("E26").Value = 16.815
$ws.Range("A27").Value = -21.757
$ws.Range("D27").Value = -8.303000000000001
$ws.Range("A28").Value = -21.875
$ws.Range("D28").Value = -7.881
$ws.Range("D30").Value = -7.157999999999999
$ws.Range("A31").Value = -21.838
$ws.Range("C31").Value = -12.849
$ws.Range("A32").Value = -21.678
$ws.Range("D32").Value = -7.822000000000001
$ws.Range("C35").Value = -12.449
$ws.Range("A36").Value = -20.32
$ws.Range("E36").Value = 16.795
$ws.Range("C37").Value = -13.76
$ws.Range("A38").Value = -19.696
$ws.Range("C39").Value = -12.745
$ws.Range("C40").Value = -12.989
$ws.Range("E42").Value = 16.547
$ws.Range("D44").Value = -7.904999999999999
$ws.Range("C45").Value = -12.569
$ws.Range("A46").Value = -21.883
$ws.Range("D47").Value = -7.531000000000001
$ws.Range("C52").Value = -11.363
$ws.Range("A54").Value = -22.15
$ws.Range("A55").Value = -22.223
$ws.Range("A56").Value = -21.997
$ws.Range("C57").Value = -13.337
$ws.Range("D58").Value = -8.218
$ws.Range("D62").Value = -7.994
$ws.Range("E64").Value = 16.979
$ws.Range("A67").Value = -21.534
$ws.Range("E68").Value = 17.26
$ws.Range("A69").Value = -21.544
$ws.Range("A72").Value = -21.481
$ws.Range("A73").Value = -19.973
$ws.Range("D77").Value = -7.949000000000001
$ws.Range("D78").Value = -7.972999999999999
$ws.Range("E79").Value = 17.266
$ws.Range("C81").Value = -12.926
$ws.Range("A83").Value = -21.628
$ws.Range("C83").Value = -12.893
$ws.Range("D84").Value = -8.019
$ws.Range("A86").Value = -22.264
$ws.Range("D89").Value = -7.233999999999999
$ws.Range("E89").Value = 17.104
$ws.Range("A91").Value = -21.595
$ws.Range("D91").Value = -6.953
$ws.Range("D92").Value = -7.007
$ws.Range("A93").Value = -21.547
$ws.Range("D96").Value = -7.423999999999999
$ws.Range("A99").Value = -19.854
$ws.Range("C100").Value = -12.215
$ws.Range("C102").Value = -13.25
$ws.Range("D102").Value = -7.602000000000001
